{"js": "// Replace the 25 division-problem answers in the single 20-row/5-column\n// table. Only every 4th row (0, 4, 8, 12, 16) holds problem text; the rows\n// in between are blank answer rows. Because some source strings repeat\n// (e.g. \"30\u00f76=5, 0\" appears twice with different replacements), cells are\n// addressed by (row, col) position rather than by text search/replace.\nconst replacements = [\n  { row: 0, col: 0, text: \"36\u00f75=7, 1\" },\n  { row: 0, col: 1, text: \"13\u00f77=1, 6\" },\n  { row: 0, col: 2, text: \"73\u00f73=24, 1\" },\n  { row: 0, col: 3, text: \"76\u00f72=38, 0\" },\n  { row: 0, col: 4, text: \"64\u00f79=7, 1\" },\n\n  { row: 4, col: 0, text: \"52\u00f73=17, 1\" },\n  { row: 4, col: 1, text: \"76\u00f73=25, 1\" },\n  { row: 4, col: 2, text: \"57\u00f77=8, 1\" },\n  { row: 4, col: 3, text: \"67\u00f75=13, 2\" },\n  { row: 4, col: 4, text: \"26\u00f77=3, 5\" },\n\n  { row: 8, col: 0, text: \"52\u00f72=26, 0\" },\n  { row: 8, col: 1, text: \"32\u00f78=4, 0\" },\n  { row: 8, col: 2, text: \"80\u00f74=20, 0\" },\n  { row: 8, col: 3, text: \"19\u00f75=3, 4\" },\n  { row: 8, col: 4, text: \"90\u00f77=12, 6\" },\n\n  { row: 12, col: 0, text: \"32\u00f72=16, 0\" },\n  { row: 12, col: 1, text: \"12\u00f73=4, 0\" },\n  { row: 12, col: 2, text: \"99\u00f79=11, 0\" },\n  { row: 12, col: 3, text: \"91\u00f72=45, 1\" },\n  { row: 12, col: 4, text: \"97\u00f79=10, 7\" },\n\n  { row: 16, col: 0, text: \"57\u00f79=6, 3\" },\n  { row: 16, col: 1, text: \"70\u00f76=11, 4\" },\n  { row: 16, col: 2, text: \"22\u00f75=4, 2\" },\n  { row: 16, col: 3, text: \"86\u00f72=43, 0\" },\n  { row: 16, col: 4, text: \"45\u00f76=7, 3\" },\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (const { row, col, text } of replacements) {\n  const cell = table.getCell(row, col);\n  cell.value = text;\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem answers in the single 20-row/5-column\n# table. Only every 4th row (rows 1, 5, 9, 13, 17 in 1-based COM indexing)\n# holds problem text; the rows in between are blank answer rows. Because\n# some source strings repeat (e.g. \"30\u00f76=5, 0\" appears twice with different\n# replacements), cells are addressed by (row, col) position rather than by\n# text search/replace.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1,1).Range.Text = \"36\u00f75=7, 1\"\n$t.Cell(1,2).Range.Text = \"13\u00f77=1, 6\"\n$t.Cell(1,3).Range.Text = \"73\u00f73=24, 1\"\n$t.Cell(1,4).Range.Text = \"76\u00f72=38, 0\"\n$t.Cell(1,5).Range.Text = \"64\u00f79=7, 1\"\n\n$t.Cell(5,1).Range.Text = \"52\u00f73=17, 1\"\n$t.Cell(5,2).Range.Text = \"76\u00f73=25, 1\"\n$t.Cell(5,3).Range.Text = \"57\u00f77=8, 1\"\n$t.Cell(5,4).Range.Text = \"67\u00f75=13, 2\"\n$t.Cell(5,5).Range.Text = \"26\u00f77=3, 5\"\n\n$t.Cell(9,1).Range.Text = \"52\u00f72=26, 0\"\n$t.Cell(9,2).Range.Text = \"32\u00f78=4, 0\"\n$t.Cell(9,3).Range.Text = \"80\u00f74=20, 0\"\n$t.Cell(9,4).Range.Text = \"19\u00f75=3, 4\"\n$t.Cell(9,5).Range.Text = \"90\u00f77=12, 6\"\n\n$t.Cell(13,1).Range.Text = \"32\u00f72=16, 0\"\n$t.Cell(13,2).Range.Text = \"12\u00f73=4, 0\"\n$t.Cell(13,3).Range.Text = \"99\u00f79=11, 0\"\n$t.Cell(13,4).Range.Text = \"91\u00f72=45, 1\"\n$t.Cell(13,5).Range.Text = \"97\u00f79=10, 7\"\n\n$t.Cell(17,1).Range.Text = \"57\u00f79=6, 3\"\n$t.Cell(17,2).Range.Text = \"70\u00f76=11, 4\"\n$t.Cell(17,3).Range.Text = \"22\u00f75=4, 2\"\n$t.Cell(17,4).Range.Text = \"86\u00f72=43, 0\"\n$t.Cell(17,5).Range.Text = \"45\u00f76=7, 3\"\n"}
